$p = $ppt.ActivePresentation
$p.Slides.Item(52).Delete()
